# Commit: "key api report excel write"
# Fills in the previously-blank rows 3-8 of the "t_接" sheet (2nd tab) with
# key-extraction API test rows, mirroring the pattern already used on the
# "t_接口py" sheet, and updates the active selections left behind on the
# "t_接" and "t_接口py" sheets.

$wb = $excel.ActiveWorkbook

# --- "t_接" sheet (2nd tab) -------------------------------------------------
$ws = $wb.Worksheets.Item(2)

# Row 3 (num=2) already had exec/title set; fill in the remaining columns.
$ws.Range("D3").Value = "save2dict"
$ws.Range("E3").Value = "auser"
$ws.Range("F3").Value = "hh"
$ws.Range("G3").Value = '${cc}mm'

# Row 4 (num=3)
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "y"
$ws.Range("C4").Value = "充3"
$ws.Range("D4").Value = "saveparam"
$ws.Range("E4").Value = "xcode"
$ws.Range("F4").Value = 20001

# Row 5 (num=4)
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "y"
$ws.Range("C5").Value = "充4"
$ws.Range("D5").Value = "get_api"
$ws.Range("E5").Value = "get"
$ws.Range("F5").Value = "za={'a':'a1'}"

# Row 6 (num=5)
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "y"
$ws.Range("C6").Value = "充5"
$ws.Range("D6").Value = "post_api"
$ws.Range("E6").Value = "post"
$ws.Range("F6").Value = 'zss=${auser}'

# Row 7 (num=7)
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "y"
$ws.Range("C7").Value = "充6"
$ws.Range("D7").Value = "get_api"
$ws.Range("E7").Value = 'get/${auser}'

# Row 8 (num=8)
$ws.Range("A8").Value = 8
$ws.Range("B8").Value = "y"
$ws.Range("C8").Value = "充7"
$ws.Range("D8").Value = "savejson"
$ws.Range("E8").Value = "next"
$ws.Range("F8").Value = "headers,X-Amzn-Trace-Id"

# Selection left on this sheet after the edit.
$ws.Range("D7").Select()

# --- "t_接口py" sheet (3rd tab) ---------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A11:XFD17").Select()
$wb.Worksheets.Item(2).Activate()
